$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date header (shared string used by B1).
# B1 holds a plain text date string ("2019-06-16" -> "2019-06-17"), but Excel's
# smart entry would otherwise auto-convert a date-like string into a date
# serial number + date style. Force text recognition, assign, then clear the
# formatting back off (restores cell to default style, keeping the text type).
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2019-06-17"
$ws.Range("B1").ClearFormats()

# Update the log level counts in column B
$ws.Range("B4").Value = 1.0
$ws.Range("B5").Value = 3.0
$ws.Range("B6").Value = 2.0
$ws.Range("B7").Value = 1.0
$ws.Range("B8").Value = 3.0
